# Horarios actualizados Línea 141 - 1179
# Update scrape timestamp across the three sheets and refresh the
# schedule rows to the newly scraped values.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "LP1912": refresh timestamp + arrival rows (row count stays the same)
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 02:13:42"

$ws1.Range("A6").Value = "02:13:42"
$ws1.Range("B6").Value = "03:09"
$ws1.Range("D6").Value = 56

$ws1.Range("A7").Value = "02:13:42"
$ws1.Range("B7").Value = "03:55"
$ws1.Range("D7").Value = 102

$ws1.Range("A8").Value = "02:13:42"
$ws1.Range("B8").Value = "04:01"
$ws1.Range("C8").Value = "81_EL PELIGRO"
$ws1.Range("D8").Value = 108

# ---------------------------------------------------------------
# Sheet "LP1912-215": refresh timestamp, update row 6, drop row 7
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 02:13:42"
$ws2.Range("A3").Value = "Total filas: 1"

$ws2.Range("A6").Value = "02:13:42"
$ws2.Range("B6").Value = "03:09"
$ws2.Range("D6").Value = 56

$ws2.Rows("7:7").Delete()

# ---------------------------------------------------------------
# Sheet "6203-6173": refresh timestamp only
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 02:13:42"
